$wb = $excel.ActiveWorkbook

# --- Sheet "Sockeye FSC and Demo": set Lake Babine Nation demo catch total for week 35 ---
$ws1 = $wb.Worksheets.Item("Sockeye FSC and Demo")
$ws1.Range("J16").Value = 34030
$ws1.Range("J16").Select()

# --- Sheet "demo catches": append 6 new Lake Babine / Gillnet rows ---
$ws3 = $wb.Worksheets.Item("demo catches")

$newRows = @(
    @(45529, 2100),
    @(45531, 840),
    @(45532, 6510),
    @(45533, 11130),
    @(45534, 5250),
    @(45535, 8200)
)

$startRow = 41
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $date = $newRows[$i][0]
    $pieces = $newRows[$i][1]

    $ws3.Cells.Item($r, 1).Value = $date
    $ws3.Cells.Item($r, 2).Value = "Lake Babine"
    $ws3.Cells.Item($r, 3).Value = "Gillnet"
    $ws3.Cells.Item($r, 4).Value = $pieces
}

# Match the date-format style used by the rest of column A (row 40 as template)
$ws3.Cells.Item(40, 1).Copy()
$ws3.Range("A41:A46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("E46").Select()
